# 2020president.xlsx tidy-up
# - Normalize the three presidential-ticket header labels (drop the stray
#   space after the underscore): "Biden_ Harris" -> "Biden_Harris", etc.
# - Remove the stray blank row that was sitting between the county data
#   and the "Check"/"Totals:" summary rows, so those two rows move up
#   one (66/67 instead of 67/68).
# - Re-apply the user's sheet view (zoomed in, column A selected) to match
#   the saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header labels for the three presidential tickets (B1:D1).
$ws.Range("B1").Value = "Biden_Harris"
$ws.Range("C1").Value = "Trump_Pence"
$ws.Range("D1").Value = "Blankenship_Mohr"

# Remove the empty row 66 that separated the county rows (2-65) from the
# "Check" formula row and the "Totals:" row, shifting those two rows up.
$ws.Rows.Item(66).Delete()

# Match the saved view: zoomed to 190% with column A selected.
$ws.Columns.Item(1).Select() | Out-Null
$excel.ActiveWindow.Zoom = 190
